# Se modifican ControlActivationActions para boton confirmar
$wb = $excel.ActiveWorkbook

# --- Rename sheet "Semilla 3" -> "Semilla 13" ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Semilla 13"

# --- Update environment config values on the renamed sheet ---
$ws2.Range("A4").Value = "10.69.60.227"
$ws2.Range("A5").Value = "10.69.60.223"
$ws2.Range("B5").Value = "DESEPOS"
$ws2.Range("A6").Value = "10.69.60.223"
$ws2.Range("A7").Value = "10.65.32.74"
$ws2.Range("B7").Value = "SIEBCBS1"

# --- Update selections / active views to match the saved workbook state ---
$ws2.Range("B20").Select() | Out-Null

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A13").Select() | Out-Null

$ws5 = $wb.Worksheets.Item(5)
$ws5.Activate() | Out-Null
$ws5.Range("A12").Select() | Out-Null
